$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = "Asuka Ikeda"
$ws.Range("J2").Value = "bdot20240415_141954/1.png"
$ws.Range("K2").Value = "「スタート」ボタンをクリックする"

# --- Row 3 ---
$ws.Range("C3").Value = "Asuka Ikeda"
$ws.Range("J3").Value = "bdot20240415_141954/2.png"
$ws.Range("K3").Value = "メニューから「設定」アイコンをクリックする"

# --- Row 4 ---
$ws.Range("C4").Value = "Asuka Ikeda"
$ws.Range("J4").Value = "bdot20240415_141954/3.png"
$ws.Range("K4").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"

# --- Row 5 ---
$ws.Range("B5").Value = "error"
$ws.Range("C5").Value = "Asuka Ikeda"
$ws.Range("J5").Value = "bdot20240415_141954/4.png"
$ws.Range("K5").Value = "0x80240fff エラー"
$ws.Range("L5").Value = "Error W"
$ws.Range("M5").Value = " エラーの Windows"

# --- Row 6 ---
$ws.Range("C6").Value = "Asuka Ikeda"
$ws.Range("J6").Value = "bdot20240415_141954/5.png"
$ws.Range("K6").Value = "デスクトップ画面の左下にある「スタート」ボタンを右クリックする"

# --- Row 7 ---
$ws.Range("B7").Value = "operation"
$ws.Range("C7").Value = "Asuka Ikeda"
$ws.Range("J7").Value = "bdot20240415_141954/5.png"
$ws.Range("K7").Value = "メニューからターミナル(管理者)をクリックする"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""

# --- Row 8 ---
$ws.Range("C8").Value = "Asuka Ikeda"
$ws.Range("J8").Value = "bdot20240415_141954/6.png"
$ws.Range("K8").Value = "ユーザーアカウント制御と表示されているウィンドウが開いたことを確認する"

# --- Row 9 ---
$ws.Range("C9").Value = "Asuka Ikeda"
$ws.Range("J9").Value = "bdot20240415_141954/7.png"
$ws.Range("K9").Value = "PowerShellウィンドウに start-transcript と入力し、[Enter]キーを押す"

# --- Row 10 ---
$ws.Range("C10").Value = "Asuka Ikeda"
$ws.Range("J10").Value = "bdot20240415_141954/8.png"
$ws.Range("K10").Value = "wuauclt.exe /resetauthorization /detectnow と入力し、[Enter]キーを押す"

# --- Row 11 ---
$ws.Range("C11").Value = "Asuka Ikeda"
$ws.Range("J11").Value = "bdot20240415_141954/9.png"
$ws.Range("K11").Value = "netsh winhttp show proxy と入力し、[Enter]キーを押す"

# --- Row 12 ---
$ws.Range("C12").Value = "Asuka Ikeda"
$ws.Range("J12").Value = "bdot20240415_141954/10.png"
$ws.Range("K12").Value = "netsh winhttp reset proxy と入力し、[Enter]キーを押す"

# --- Row 13 ---
$ws.Range("C13").Value = "Asuka Ikeda"
$ws.Range("J13").Value = "bdot20240415_141954/1.png"
$ws.Range("K13").Value = "「スタート」ボタンをクリックする"

# --- Row 14 ---
$ws.Range("C14").Value = "Asuka Ikeda"
$ws.Range("J14").Value = "bdot20240415_141954/2.png"
$ws.Range("K14").Value = "メニューから「設定」アイコンをクリックする"

# --- Row 15 ---
$ws.Range("C15").Value = "Asuka Ikeda"
$ws.Range("J15").Value = "bdot20240415_141954/3.png"
$ws.Range("K15").Value = "左側のメニューからWindows Updateをクリックし、Windows Update画面に移動する"

# --- Row 16 ---
$ws.Range("C16").Value = "Asuka Ikeda"
$ws.Range("J16").Value = "bdot20240415_141954/11.png"
$ws.Range("K16").Value = "「更新プログラムのチェック」ボタンをクリックする"

Write-Output "done"
